# Update Notification Center: default to tomorrow, add send button, and refresh excel report
# This script refreshes the "reporte_agendas" worksheet data:
#  - Marks row 63 (RT 3D session on 10/12/2025) as COMPLETADO
#  - Appends 5 new appointment rows (74-78) at the bottom of the report

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark existing pending row as completed
$ws.Range("L63").Value = "COMPLETADO"

# New rows to append to the report
$newRows = @(
    @("10/12/2025", "18:00:00", "53166674", "AGUIRRE, GUILLERMINA", "3794630005", "IOSCOR", "F", 12, "EEG", "ELECTRO Y MAPEOS", "DR. MONZON ROMILIO", "COMPLETADO"),
    @("17/12/2025", "18:00:00", "30141716", "AGUIRRE, ALFREDO RAUL", "3794630005", "OSDE", "M", 42, "TAC COMPLETA DE ABDOMEN", "TOMOGRAFIAS Y RX", "DR. MONZON ROMILIO", "COMPLETADO"),
    @("17/12/2025", "19:00:00", "30141716", "AGUIRRE, ALFREDO RAUL", "3794630005", "OSDE", "M", 42, "TAC DE CEREBRO", "TOMOGRAFIAS Y RX", "RAINERO FEDERICO", "COMPLETADO"),
    @("17/12/2025", "19:00:00", "30141716", "AGUIRRE, ALFREDO RAUL", "3794630005", "OSDE", "M", 42, "TAC DE CUELLO", "TOMOGRAFIAS Y RX", "RAINERO FEDERICO", "COMPLETADO"),
    @("22/12/2025", "08:00:00", "31648199", "VALLEJOS SCHULZE, MARIA ELENA", "3794774785", "IOSCOR", "F", 40, "CENTELLOGRAMA DE TIROIDES", "CAMARA GAMMA", "DE LOS REYES", "PENDIENTE")
)

$startRow = 74
$endRow = $startRow + $newRows.Count - 1

# Format all the new cells as Text first (except column H, which is numeric),
# so Excel does not auto-convert date-like / numeric-looking strings
# (e.g. "10/12/2025", "30141716") into dates or numbers.
$ws.Range("A$startRow`:G$endRow").NumberFormat = "@"
$ws.Range("I$startRow`:L$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
}

# Drop the temporary Text formatting now that the values are locked in as
# text, so the new rows end up with the same default (unstyled) look as the
# rest of the report.
$ws.Range("A$startRow`:L$endRow").ClearFormats()
